# Actualizacion Datos Personales 4 nov
$wb = $excel.ActiveWorkbook

$ws5AEM = $wb.Worksheets.Item("5AEM")
$ws5BEM = $wb.Worksheets.Item("5BEM")

# --- Corrections to existing values (text fixes) ---

# 5AEM: correct student's own e-mail typo (row 18, Yair Antonio Perez Romero)
$ws5AEM.Range("E18").Value = "yair26prz@gmail.com"

# 5BEM: correct a phone number (row 13, Juan Antonio Flores de la Cruz)
$ws5BEM.Range("F13").Value = "6631083766"

# 5BEM: replace placeholder tutor name with the real tutor (row 8, Raul Arturo Citlahua Hernandez)
$ws5BEM.Range("H8").Value = "MARÍA FILOMENA HERNÁNDEZ CHONCOA"

# 5BEM: correct tutor name (row 25, Arian Alexis Muñoz Luna)
$ws5BEM.Range("H25").Value = "JUAN MAURICIO MUÑOZ MARTINEZ"

# --- New tutor contact information added ---

# 5AEM row 7: Cristian Antonio Carrasco Sandoval
$ws5AEM.Range("H7").Value = "BLANCA ESTELA SANDOVAL DÍAZ"
$ws5AEM.Range("I7").Value = "as6198099@gmail.com"
$ws5AEM.Range("J7").NumberFormat = "@"
$ws5AEM.Range("J7").Value = "2721417437"

# 5AEM row 37: Emilio Valderrama Rodriguez
$ws5AEM.Range("H37").Value = "MARÍA TEREZA RÓDRIGUEZ LOPEZ"
$ws5AEM.Range("I37").Value = "Maytequila133@gmail.com"
$ws5AEM.Range("J37").NumberFormat = "@"
$ws5AEM.Range("J37").Value = "2721270249"

# 5BEM row 8: Raul Arturo Citlahua Hernandez - tutor phone number
$ws5BEM.Range("J8").NumberFormat = "@"
$ws5BEM.Range("J8").Value = "2722480188"

# 5BEM row 13: Juan Antonio Flores de la Cruz - tutor info
$ws5BEM.Range("H13").Value = "GABRIELA FLORES DE LA CRUZ"
$ws5BEM.Range("J13").NumberFormat = "@"
$ws5BEM.Range("J13").Value = "2721135977"

$wb.Save()
